# "updated thread local driver" - refresh the test-data emails used by the
# ecommerce test suite and leave the workbook positioned on Sheet1 (the data
# sheet) instead of the login sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Roll the per-row test emails forward: drop the tom13/14/15 batch and mint
# tom18-21 in their place (tom17 -> tom18 on row 2, tom13/14/15 -> tom19/20/21
# on rows 3-5).
$ws1.Range("C2").Value = "tom18@gmail.com"
$ws1.Range("C3").Value = "tom19@gmail.com"
$ws1.Range("C4").Value = "tom20@gmail.com"
$ws1.Range("C5").Value = "tom21@gmail.com"

# Sheet1 becomes the active/selected sheet & cell (was "login" before).
$ws1.Activate() | Out-Null
$ws1.Range("D8").Select() | Out-Null
